# Apply updated view-count / price values to "展览" (sheet 1) and "全部类型" (sheet 4) worksheets
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3139
$wsExhibit.Range("F4").Value = 1096
$wsExhibit.Range("F5").Value = 86
$wsExhibit.Range("G5").Value = 25
$wsExhibit.Range("F6").Value = 37
$wsExhibit.Range("F7").Value = 274
$wsExhibit.Range("F9").Value = 1129
$wsExhibit.Range("F10").Value = 15764
$wsExhibit.Range("F11").Value = 244
$wsExhibit.Range("F12").Value = 184
$wsExhibit.Range("F13").Value = 1027
$wsExhibit.Range("F14").Value = 6190
$wsExhibit.Range("F15").Value = 626
$wsExhibit.Range("F18").Value = 9
$wsExhibit.Range("F19").Value = 119
$wsExhibit.Range("F22").Value = 634
$wsExhibit.Range("F23").Value = 17
$wsExhibit.Range("F25").Value = 3
$wsExhibit.Range("F27").Value = 867
$wsExhibit.Range("F28").Value = 30
$wsExhibit.Range("F29").Value = 5003
$wsExhibit.Range("F30").Value = 490
$wsExhibit.Range("F31").Value = 11082
$wsExhibit.Range("F32").Value = 1232
$wsExhibit.Range("F33").Value = 14
$wsExhibit.Range("F34").Value = 124
$wsExhibit.Range("F35").Value = 170
$wsExhibit.Range("F36").Value = 3806
$wsExhibit.Range("F37").Value = 265

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3139
$wsAll.Range("F5").Value = 1096
$wsAll.Range("F6").Value = 86
$wsAll.Range("G6").Value = 25
$wsAll.Range("F7").Value = 37
$wsAll.Range("F8").Value = 274
$wsAll.Range("F10").Value = 1129
$wsAll.Range("F11").Value = 15764
$wsAll.Range("F12").Value = 244
$wsAll.Range("F13").Value = 184
$wsAll.Range("F14").Value = 1027
$wsAll.Range("F15").Value = 6190
$wsAll.Range("F16").Value = 626
$wsAll.Range("F19").Value = 9
$wsAll.Range("F20").Value = 119
$wsAll.Range("F23").Value = 634
$wsAll.Range("F24").Value = 17
$wsAll.Range("F26").Value = 3
$wsAll.Range("F28").Value = 867
$wsAll.Range("F29").Value = 30
$wsAll.Range("F30").Value = 5003
$wsAll.Range("F31").Value = 490
$wsAll.Range("F33").Value = 11082
$wsAll.Range("F34").Value = 1232
$wsAll.Range("F35").Value = 14
$wsAll.Range("F36").Value = 124
$wsAll.Range("F37").Value = 170
$wsAll.Range("F38").Value = 3806
$wsAll.Range("F39").Value = 265

